$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Items")

# --- Header row (row 2): "Effect" / "Purpose of the quiz question" ---
$ws.Range("B2").Value = "Effect"
$ws.Range("B2").HorizontalAlignment = -4108
$ws.Range("C2").Value = "Purpose of the quiz question"

# --- Fill in Effect / Purpose columns for the existing rows ---
$ws.Range("B4").Value = "random, present-tense verb"
$ws.Range("C4").Value = "Student should translate the word into Japanese"

$ws.Range("B5").Value = "random English word"
$ws.Range("C5").Value = "Student should translate the word into Japanese"

$ws.Range("B6").Value = "random, difficult to pronounce English word"
$ws.Range("C6").Value = "Student should pronounce the word properly."

# Row 7 (SpinyBeetle) stays without Effect/Purpose.

# --- Insert a new item row for PirahnaPlant before the PowButton row ---
$ws.Rows.Item(8).Insert()
$ws.Range("A8").Value = "PirahnaPlant"
$ws.Range("B8").Value = "random, English target sentence"
$ws.Range("C8").Value = "Student should translate the target sentence into Japanese."

# PowButton (now row 9) gets its Effect / Purpose values
$ws.Range("B9").Value = "random image"
$ws.Range("C9").Value = "Student should identify the thing(s) in the image."

# Bombomb (now row 10) gets its Effect / Purpose values
$ws.Range("B10").Value = "randomly mix items on the screen"
$ws.Range("C10").Value = "N/A"

# FireFlower (row 11) and BlueCoin (row 12) stay without Effect/Purpose.

# --- Append a new item row for QuestionBlock ---
$ws.Range("A13").Value = "QuestionBlock"
$ws.Range("A13").HorizontalAlignment = -4131
$ws.Range("A13").IndentLevel = 1
$ws.Range("B13").Value = "shows ""Ask me a question."""
$ws.Range("C13").Value = "Student should ask a question in English."

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 14.166666666666666
$ws.Columns.Item(2).ColumnWidth = 57
$ws.Columns.Item(3).ColumnWidth = 58.666666666666664

# --- Selection matches the saved workbook state ---
$ws.Range("B12").Select()
